# Apply the "Generate Report for handoff" update:
#  - A new localization handoff was produced for source file
#    306b57c1-daec-45e0-b492-72653faf7d85.md (previously ae3af397-...md)
#  - A second source file, f2bb19fd-6824-4ed0-8680-0628834afc01.md, failed
#    its handoff transform, so a new status row is appended after the
#    existing ".localization-config" row (which keeps its old content but
#    moves down to the new last row).

$wb = $excel.ActiveWorkbook

$newMd        = "306b57c1-daec-45e0-b492-72653faf7d85.md"
$newConfigMd  = "f2bb19fd-6824-4ed0-8680-0628834afc01.md"
$configFile   = ".localization-config"

$newXlfZh = "306b57c1-daec-45e0-b492-72653faf7d85.4752fa7747eed55319e669b57eba3ac6d8468e58.zh-cn.xlf"
$newXlfDe = "306b57c1-daec-45e0-b492-72653faf7d85.4752fa7747eed55319e669b57eba3ac6d8468e58.de-de.xlf"

$readyStatus  = "Ready for handoff"
$failedStatus = "Handoff transform failed"
$notLocalized = "Not to be localized"

$handoffDateZh = "2016-01-26 04:54:25"
$handoffDateDe = "2016-01-26 04:54:34"
$zeroDate      = "0001-01-01 00:00:00"
$includeReason = "Include"
$ignoredReason = "Ignored"

$urlRepoRoot = "https://github.com/OpenLocalizationTest/oltest/blob/baf040cf35a7c8054bf801d85404441a1487b81d"
$urlZhHandoff = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/047cba0713f6b859a3a72a2157df7f10a3bf5b9e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
$urlDeHandoff = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6a0e7f458f02885d508fa4f8db177423fa0835f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho"

function Set-FileHyperlink($ws, $cellAddr, $displayText, $url) {
    $ws.Range($cellAddr).Value = $displayText
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
    # Match the workbook's existing hyperlink look (underlined, cornflower-blue text)
    # rather than the engine's default "Hyperlink" cell style.
    $ws.Range($cellAddr).Font.Underline = 2
    $ws.Range($cellAddr).Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

Set-FileHyperlink $wsOverview "A2" $newMd "$urlRepoRoot/e2e/$newMd"
$wsOverview.Range("B2").Value = $readyStatus
$wsOverview.Range("C2").Value = $readyStatus

Set-FileHyperlink $wsOverview "A3" $newConfigMd "$urlRepoRoot/e2e/$newConfigMd"
$wsOverview.Range("B3").Value = $failedStatus
$wsOverview.Range("C3").Value = $failedStatus

Set-FileHyperlink $wsOverview "A4" $configFile "$urlRepoRoot/$configFile"
$wsOverview.Range("B4").Value = $notLocalized
$wsOverview.Range("C4").Value = $notLocalized

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

Set-FileHyperlink $wsZh "A2" $newMd "$urlRepoRoot/e2e/$newMd"
$wsZh.Range("B2").Value = $readyStatus
Set-FileHyperlink $wsZh "C2" $newXlfZh "$urlZhHandoff/$newXlfZh"
$wsZh.Range("D2").Value = $handoffDateZh
$wsZh.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G2").Value = $zeroDate
$wsZh.Range("H2").Value = $includeReason

Set-FileHyperlink $wsZh "A3" $newConfigMd "$urlRepoRoot/e2e/$newConfigMd"
$wsZh.Range("B3").Value = $failedStatus
$wsZh.Range("D3").Value = $zeroDate
$wsZh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G3").Value = $zeroDate
$wsZh.Range("H3").Value = $ignoredReason

Set-FileHyperlink $wsZh "A4" $configFile "$urlRepoRoot/$configFile"
$wsZh.Range("B4").Value = $notLocalized
$wsZh.Range("D4").Value = $zeroDate
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G4").Value = $zeroDate
$wsZh.Range("H4").Value = $ignoredReason

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

Set-FileHyperlink $wsDe "A2" $newMd "$urlRepoRoot/e2e/$newMd"
$wsDe.Range("B2").Value = $readyStatus
Set-FileHyperlink $wsDe "C2" $newXlfDe "$urlDeHandoff/$newXlfDe"
$wsDe.Range("D2").Value = $handoffDateDe
$wsDe.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G2").Value = $zeroDate
$wsDe.Range("H2").Value = $includeReason

Set-FileHyperlink $wsDe "A3" $newConfigMd "$urlRepoRoot/e2e/$newConfigMd"
$wsDe.Range("B3").Value = $failedStatus
$wsDe.Range("D3").Value = $zeroDate
$wsDe.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G3").Value = $zeroDate
$wsDe.Range("H3").Value = $ignoredReason

Set-FileHyperlink $wsDe "A4" $configFile "$urlRepoRoot/$configFile"
$wsDe.Range("B4").Value = $notLocalized
$wsDe.Range("D4").Value = $zeroDate
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G4").Value = $zeroDate
$wsDe.Range("H4").Value = $ignoredReason
